# EMLtools datatablename_attributes_draft.xlsx header rework
#
# The template's single header row moves from the 4-column
# "description / columnName / unitOrCodeExplanationOrDateFormat / emptyValueCode"
# layout to the 7-column EML `attributeName` layout:
#   attributeName | attributeDefinition | class | unit |
#   dateTimeFormatString | missingValueCode | missingValueCodeExplanation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (A1:G1) -- this also rebuilds the shared-string table
# in the same order the labels are written here.
$ws.Range("A1").Value = "attributeName"
$ws.Range("B1").Value = "attributeDefinition"
$ws.Range("C1").Value = "class"
$ws.Range("D1").Value = "unit"
$ws.Range("E1").Value = "dateTimeFormatString"
$ws.Range("F1").Value = "missingValueCode"
$ws.Range("G1").Value = "missingValueCodeExplanation"

# Column widths matching the published template (best-fit-style widths,
# expressed in ColumnWidth "characters" units).
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 17.333333333333336
$ws.Columns.Item(3).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws.Columns.Item(5).ColumnWidth = 20.333333333333336
$ws.Columns.Item(6).ColumnWidth = 16.833333333333336
$ws.Columns.Item(7).ColumnWidth = 27.666666666666668
